$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.337.22"
$ws.Range("E2").Value = "  -0.12%  "

$ws.Range("D3").Value = "1.879.68"
$ws.Range("E3").Value = "  +0.29%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "'0.7100"
$ws.Range("E5").Value = "  -0.22%  "

$ws.Range("D6").Value = "'242.42"
$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").Value = "'0.08025"
$ws.Range("E8").Value = "  +3.15%  "

$ws.Range("D9").Value = "'0.3139"
$ws.Range("E9").Value = "  +0.75%  "

$ws.Range("D10").Value = "'25.11"
$ws.Range("E10").Value = "  +0.40%  "

$ws.Range("D11").Value = "'0.08328"
$ws.Range("E11").Value = "  -1.65%  "

$ws.Range("D12").Value = "1.899.23"
$ws.Range("E12").Value = "  +1.47%  "

$ws.Range("D13").Value = "'5.270"
$ws.Range("E13").Value = "  +0.52%  "

$ws.Range("D14").Value = "'94.63"
$ws.Range("E14").Value = "  +3.74%  "

$ws.Range("D15").Value = "'0.7179"
$ws.Range("E15").Value = "  +0.76%  "

$ws.Range("D16").Value = "'6.356"
$ws.Range("E16").Value = "  +5.28%  "

$ws.Range("D17").Value = "'0.000008685"
$ws.Range("E17").Value = "  +5.37%  "

$ws.Range("D18").Value = "29.346.15"
$ws.Range("E18").Value = "  -0.11%  "

$ws.Range("D19").Value = "'243.07"
$ws.Range("E19").Value = "  +1.00%  "

$ws.Range("D20").Value = "2.141.42"
$ws.Range("E20").Value = "  +1.06%  "

$ws.Range("E21").Value = "  +0.39%  "

$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("D23").Value = "'7.862"
$ws.Range("E23").Value = "  +0.80%  "

$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("D25").Value = "'0.1576"
$ws.Range("E25").Value = "  -1.91%  "

$ws.Range("D26").Value = "'163.61"
$ws.Range("E26").Value = "  -0.09%  "

$ws.Range("D27").Value = "'9.075"
$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("E28").Value = "  +0.64%  "

$ws.Range("D29").Value = "'1.507"
$ws.Range("E29").Value = "  -0.21%  "

$ws.Range("D30").Value = "'4.437"
$ws.Range("E30").Value = "  +0.27%  "

$ws.Range("D31").Value = "'4.366"
$ws.Range("E31").Value = "  +1.11%  "

$ws.Range("D32").Value = "'1.203"
$ws.Range("E32").Value = "  -6.04%  "

$ws.Range("D33").Value = "'0.05399"
$ws.Range("E33").Value = "  +1.83%  "

$ws.Range("D34").Value = "'1.940"
$ws.Range("E34").Value = "  +0.19%  "

$ws.Range("D35").Value = "'0.7756"
$ws.Range("E35").Value = "  +3.85%  "

$ws.Range("D36").Value = "'1.178"
$ws.Range("E36").Value = "  +0.07%  "

$ws.Range("D37").Value = "'2.686"
$ws.Range("E37").Value = "  -0.38%  "

$ws.Range("D38").Value = "'0.01885"
$ws.Range("E38").Value = "  +0.69%  "

$ws.Range("D39").Value = "1.270.78"
$ws.Range("E39").Value = "  +5.30%  "

$ws.Range("D40").Value = "'2.745"
$ws.Range("E40").Value = "  +0.90%  "

$ws.Range("D41").Value = "'6.550"

$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "'113.77"
$ws.Range("E42").Value = "  +5.14%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.9186"
$ws.Range("E43").Value = "  +3.44%  "

$ws.Range("D44").Value = "'74.62"
$ws.Range("E44").Value = "  +2.34%  "

$ws.Range("D45").Value = "'1.000"
$ws.Range("E45").Value = "  +0.05%  "

$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "2.049.99"
$ws.Range("E46").Value = "  +1.47%  "

$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "'0.00000000127"
$ws.Range("E47").Value = "  +3.89%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.5223"
$ws.Range("E48").Value = "  +0.27%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'1.805"
$ws.Range("E49").Value = "  -0.81%  "

$ws.Range("D50").Value = "'9.565"
$ws.Range("E50").Value = "  +1.98%  "

$ws.Range("E51").Value = "  +1.26%  "
